# Add a new user row ("José" / "rex") to the Usuarios sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")

$ws.Range("A6").Value = "José"
$ws.Range("B6").Value = "rex"

# Reuse the existing "12345678" text value (same password as row 4) via
# copy/paste so the cell keeps its original shared-string/text type
# instead of being reinterpreted as a number.
$ws.Range("C4").Copy($ws.Range("C6"))

$ws.Range("D6").Value = "cuentasfalsasxdd@hotmail.com"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 5
